$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ULBI")

$ws.Columns("D:D").Insert()

$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("D5:D6").Clear()

$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 87200
$ws.Range("D9").Value = 61600
$ws.Range("D10").Value = 25600
$ws.Range("D12").Value = 4400
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 400
$ws.Range("D17").Value = 80600
$ws.Range("D18").Value = 6600
$ws.Range("D20").Value = 100
$ws.Range("D21").Value = 9000
$ws.Range("D22").Value = 100
$ws.Range("D23").Value = 6600
$ws.Range("D24").Value = -18400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 25000
$ws.Range("D27").Value = 24900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -100
$ws.Range("D33").Value = 24900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 24900
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 25600
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 16000
$ws.Range("D44").Value = 22800
$ws.Range("D45").Value = 2800
$ws.Range("D46").Value = 67200
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 10700
$ws.Range("D49").Value = 26600
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 15500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 120100
$ws.Range("D57").Value = 9900
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 5000
$ws.Range("D60").Value = 14900
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 15500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -58000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 104600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 24900
$ws.Range("D83").Value = 2400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 10900
$ws.Range("D91").Value = -4200
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -4200
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 1200
$ws.Range("D101").Value = -300
$ws.Range("D102").Value = 7600

Write-Output "done"
